$wb = $excel.ActiveWorkbook

# --- Update the conversion note on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 8.48 = 35774.92 pesos`n✅ 35774.92 pesos = 8.47 = 956.82 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 117.96
$ws2.Range("O10").Value = 4220.01

$ws2.Range("N12").Value = 4225
$ws2.Range("O12").Value = 113
